{"js": "// Ngoc Duc pham vi du an\n// \"Ph\u1ea1m vi d\u1ef1 \u00e1n\" section: expand the italic sentence about platform\n// support from \"Ch\u1ea1y \u0111\u01b0\u1ee3c tr\u00ean n\u1ec1n t\u1ea3ng n\u00e0o, OS n\u00e0o?...\" to\n// \"Ch\u1ea1y \u0111\u01b0\u1ee3c tr\u00ean m\u1ecdi n\u1ec1n t\u1ea3ng m\u00e0 h\u1ed7 tr\u1ee3 Microsoft Temas.\", landing the\n// new wording in five separate (but identically-formatted, italic) runs\n// exactly as the author produced them interactively in Word.\n\nconst body = context.document.body;\n\n// The sentence lives in its own (italic) paragraph right after the\n// \"Ph\u1ea1m vi d\u1ef1 \u00e1n\" Heading2. Locate it via the old, full sentence so we\n// only ever touch that one paragraph even if similar words appear\n// elsewhere in the report.\nconst oldSentence =\n  \"Ch\u1ea1y \u0111\u01b0\u1ee3c tr\u00ean n\u1ec1n t\u1ea3ng n\u00e0o, OS n\u00e0o?...\";\n\nconst results = body.search(oldSentence, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\n    \"edit.js: could not find the target sentence '\" + oldSentence + \"'\"\n  );\n}\n\nconst target = results.items[0];\n\n// Build the replacement as raw OOXML (FlatOPC) so the five runs land as\n// independent <w:r> elements instead of being coalesced into one run by\n// the higher-level insertText() API. The first run keeps the original\n// w:rsidRPr=\"00A01B4E\" attribute carried over from the source run; the\n// rest are plain new runs, matching the canonical edit.\nconst flatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r w:rsidRPr=\"00A01B4E\"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Ch\u1ea1y \u0111\u01b0\u1ee3c tr\u00ean</w:t></w:r>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\"preserve\"> m\u1ecdi</w:t></w:r>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\"preserve\"> n\u1ec1n t</w:t></w:r>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>\u1ea3ng</w:t></w:r>\n            <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\"preserve\"> m\u00e0 h\u1ed7 tr\u1ee3 Microsoft Temas.</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Ngoc Duc pham vi du an\n# \"Ph\u1ea1m vi d\u1ef1 \u00e1n\" section: expand the italic sentence about platform\n# support from \"Ch\u1ea1y \u0111\u01b0\u1ee3c tr\u00ean n\u1ec1n t\u1ea3ng n\u00e0o, OS n\u00e0o?...\" to\n# \"Ch\u1ea1y \u0111\u01b0\u1ee3c tr\u00ean m\u1ecdi n\u1ec1n t\u1ea3ng m\u00e0 h\u1ed7 tr\u1ee3 Microsoft Temas.\", landing the\n# new wording in five separate (but identically-formatted, italic) runs\n# exactly as the author produced them interactively in Word.\n\n$d = $word.ActiveDocument\n\n$oldSentence = \"Ch\u1ea1y \u0111\u01b0\u1ee3c tr\u00ean n\u1ec1n t\u1ea3ng n\u00e0o, OS n\u00e0o?...\"\n\n# Scope the Find to the whole story, then re-materialise the hit as an\n# explicit Range(start, end) -- InsertXML only replaces content when\n# called on a freshly-minted Range object, not the live Find range.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = $oldSentence\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"edit.ps1: could not find the target sentence '$oldSentence'\"\n}\n\n$target = $d.Range($rng.Start, $rng.End)\n\n# Raw OOXML (FlatOPC) keeps the five runs independent (<w:r> per chunk)\n# instead of the COM bridge silently coalescing same-formatted runs the\n# way Range.Text / Range.InsertAfter would. First run keeps the original\n# w:rsidRPr=\"00A01B4E\" carried over from the source run.\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r w:rsidRPr=\"00A01B4E\"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Ch\u1ea1y \u0111\u01b0\u1ee3c tr\u00ean</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\"preserve\"> m\u1ecdi</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\"preserve\"> n\u1ec1n t</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>\u1ea3ng</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\"preserve\"> m\u00e0 h\u1ed7 tr\u1ee3 Microsoft Temas.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$target.InsertXML($flatOpc)\n"}
